$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) figures.
# D-column prices are stored as text (e.g. "591.95", "67.369.47"); force
# the Text number format before assigning so Excel does not silently
# convert number-looking strings into numeric cell values.

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '67.369.47'
$ws.Cells.Item(2, 5).Value = '  -0.53%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.614.14'
$ws.Cells.Item(3, 5).Value = '  +0.02%  '
$ws.Cells.Item(4, 5).Value = '  -0.05%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '591.95'
$ws.Cells.Item(5, 5).Value = '  -0.24%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '151.22'
$ws.Cells.Item(6, 5).Value = '  -2.49%  '
$ws.Cells.Item(7, 5).Value = '  +0.07%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.550'
$ws.Cells.Item(8, 5).Value = '  +0.38%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '2.612.94'
$ws.Cells.Item(9, 5).Value = '  -0.04%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.123'
$ws.Cells.Item(10, 5).Value = '  -2.69%  '
$ws.Cells.Item(11, 5).Value = '  +0.08%  '
$ws.Cells.Item(12, 5).Value = '  -0.73%  '
$ws.Cells.Item(13, 5).Value = '  -1.98%  '
$ws.Cells.Item(14, 5).Value = '  +0.03%  '
$ws.Cells.Item(16, 5).Value = '  -3.61%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '67.257.89'
$ws.Cells.Item(17, 5).Value = '  -0.45%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '2.616.89'
$ws.Cells.Item(18, 5).Value = '  +0.26%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '367.52'
$ws.Cells.Item(19, 5).Value = '  +1.05%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '11.06'
$ws.Cells.Item(20, 5).Value = '  -1.03%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '7.36'
$ws.Cells.Item(21, 5).Value = '  -4.00%  '
$ws.Cells.Item(22, 5).Value = '  -0.35%  '
$ws.Cells.Item(23, 5).Value = '  -1.91%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.04'
$ws.Cells.Item(24, 5).Value = '  +1.48%  '
$ws.Cells.Item(25, 5).Value = '  +0.01%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '10.08'
$ws.Cells.Item(26, 5).Value = '  +3.44%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '66.16'
$ws.Cells.Item(27, 5).Value = '  -1.77%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '582.73'
$ws.Cells.Item(29, 5).Value = '  +1.14%  '
$ws.Cells.Item(31, 5).Value = '  -3.21%  '
$ws.Cells.Item(32, 5).Value = '  -3.91%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '7.66'
$ws.Cells.Item(33, 5).Value = '  -3.50%  '
$ws.Cells.Item(34, 5).Value = '  -2.81%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.00'
$ws.Cells.Item(35, 5).Value = '  +0.15%  '
$ws.Cells.Item(36, 5).Value = '  -5.59%  '
$ws.Cells.Item(37, 5).Value = '  -1.78%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '155.48'
$ws.Cells.Item(38, 5).Value = '  -2.03%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '18.99'
$ws.Cells.Item(39, 5).Value = '  -1.72%  '
$ws.Cells.Item(40, 5).Value = '  +2.28%  '
$ws.Cells.Item(41, 5).Value = '  -1.12%  '
$ws.Cells.Item(42, 5).Value = '  -1.75%  '
$ws.Cells.Item(43, 5).Value = '  -0.70%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '16.81'
$ws.Cells.Item(44, 5).Value = '  +2.48%  '
$ws.Cells.Item(45, 5).Value = '  -0.09%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '155.47'
$ws.Cells.Item(46, 5).Value = '  -0.10%  '
$ws.Cells.Item(47, 5).Value = '  +2.02%  '
$ws.Cells.Item(48, 5).Value = '  -0.60%  '
$ws.Cells.Item(49, 5).Value = '  -1.24%  '
$ws.Cells.Item(50, 5).Value = '  +0.68%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '21.72'
$ws.Cells.Item(51, 5).Value = '  +4.25%  '
